# Update TPM-derived computed values in the Fgf2-Fgfr1 NATMI output sheet.
# The underlying values changed (new TPM input), so we overwrite the
# previously-computed numeric columns (G..J, M..T) with their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value
$updates = @{
    2  = @{ G=0.1636683333333333; H=0.491005;  I=0.008639493057305454; J=0.008639493057305455;
            M=1.845768666666667;  N=5.537306;  O=0.01459089321241885;  P=0.01459089321241885;
            Q=0.3020938813922222; R=2.71884493253; S=0.0001260579206085779; T=0.000126057920608578 }
    3  = @{ G=0.1636683333333333; H=0.491005;  I=0.008639493057305454; J=0.008639493057305455;
            O=0.6557810310272387; P=0.6557810310272387;
            Q=13.57747151749389;  R=122.197243657445; S=0.005665615664672441; T=0.005665615664672443 }
    4  = @{ G=0.1636683333333333; H=0.491005;  I=0.008639493057305454; J=0.008639493057305455;
            M=41.69841866666667;  N=125.095256; O=0.3296280757603424;  P=0.3296280757603424;
            Q=6.824710685808889;  R=61.42239617228; S=0.002847819472024434; T=0.002847819472024435 }
    5  = @{ I=0.808839719627903;  J=0.8088397196279031;
            M=1.845768666666667;  N=5.537306;  O=0.01459089321241885;  P=0.01459089321241885;
            Q=28.28239211558533;  R=254.541529040268; S=0.01180169397505354;  T=0.01180169397505354 }
    6  = @{ I=0.808839719627903;  J=0.8088397196279031;
            O=0.6557810310272387; P=0.6557810310272387;
            S=0.5304217452733689; T=0.530421745273369 }
    7  = @{ I=0.808839719627903;  J=0.8088397196279031;
            M=41.69841866666667;  N=125.095256; O=0.3296280757603424;  P=0.3296280757603424;
            Q=638.9376137044854;  R=5750.438523340369; S=0.2666162803794805; T=0.2666162803794806 }
    8  = @{ G=3.457711333333334;  H=10.373134;  I=0.1825207873147914;  J=0.1825207873147914;
            M=1.845768666666667;  N=5.537306;  O=0.01459089321241885;  P=0.01459089321241885;
            Q=6.382135237444889;  R=57.439217137004; S=0.002663141316756735; T=0.002663141316756735 }
    9  = @{ G=3.457711333333334;  H=10.373134;  I=0.1825207873147914;  J=0.1825207873147914;
            O=0.6557810310272387; P=0.6557810310272387;
            Q=286.8421532003696; R=2581.579378803326; S=0.1196936700891973;  T=0.1196936700891973 }
    10 = @{ G=3.457711333333334;  H=10.373134;  I=0.1825207873147914;  J=0.1825207873147914;
            M=41.69841866666667;  N=125.095256; O=0.3296280757603424;  P=0.3296280757603424;
            Q=144.1810948058116; R=1297.629853252304; S=0.06016397590883742; T=0.06016397590883742 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
